$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'56.891.50"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = "'2.317.41"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'529.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = "'132.15"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('D7').Value = "'0.995"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('D9').Value = "'2.343.77"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = "'5.29"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = "'2.736.69"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = "'23.47"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.69%  '
$ws.Range('D16').Value = "'56.931.23"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = "'2.328.40"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').Value = "'336.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'10.42"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'6.85"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = "'61.53"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').Value = "'8.70"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = "'172.68"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.99%  '
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').Value = "'0.0₃0724"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.53%  '
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').Value = "'18.47"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D35').Value = "'0.993"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = "'1.25"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('D37').Value = "'0.925"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('D38').Value = "'3.97"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('D39').Value = "'39.24"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('E40').Value = '  -3.37%  '
$ws.Range('D41').Value = "'5.76"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.69%  '
$ws.Range('D42').Value = "'149.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').Value = "'282.60"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').Value = "'0.0500"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('D48').Value = "'18.84"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('D49').Value = "'0.558"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('E51').Value = '  -0.92%  '
